$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "1327708"
$ws.Cells.Item(2, 2).Value = "https://aiesec.org/opportunity/global-talent/1327708"
$ws.Cells.Item(2, 3).Value = "[CSI] Finance Analyst"
$ws.Cells.Item(2, 4).Value = "Charles-de-Gaulle-Straße 20, 53113 Bonn, Germany"
$ws.Cells.Item(2, 5).Value = "Yes"
$ws.Cells.Item(2, 6).Value = "2 applicants"
$ws.Cells.Item(2, 7).Value = "6 - 18 Months"
$ws.Cells.Item(2, 8).Value = "DHL Group"
$ws.Cells.Item(2, 5).Interior.Color = 65535

$ws.Cells.Item(3, 1).Value = "1327712"
$ws.Cells.Item(3, 2).Value = "https://aiesec.org/opportunity/global-talent/1327712"
$ws.Cells.Item(3, 3).Value = "Business Process Management Intern"
$ws.Cells.Item(3, 4).Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Cells.Item(3, 5).Value = "No"
$ws.Cells.Item(3, 6).Value = "1 applicant"
$ws.Cells.Item(3, 7).Value = "6 - 18 Months"
$ws.Cells.Item(3, 8).Value = "Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"

$ws.Cells.Item(4, 1).Value = "1327689"
$ws.Cells.Item(4, 2).Value = "https://aiesec.org/opportunity/global-talent/1327689"
$ws.Cells.Item(4, 3).Value = "Commercial Analyst"
$ws.Cells.Item(4, 4).Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Cells.Item(4, 5).Value = "No"
$ws.Cells.Item(4, 6).Value = "2 applicants"
$ws.Cells.Item(4, 7).Value = "6 - 18 Months"
$ws.Cells.Item(4, 8).Value = "Organon Pharma"

$ws.Cells.Item(5, 1).Value = "1327593"
$ws.Cells.Item(5, 2).Value = "https://aiesec.org/opportunity/global-talent/1327593"
$ws.Cells.Item(5, 3).Value = "Marketing Assistant/Support"
$ws.Cells.Item(5, 4).Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Cells.Item(5, 5).Value = "No"
$ws.Cells.Item(5, 6).Value = "5 applicants"
$ws.Cells.Item(5, 7).Value = "6 - 18 Months"
$ws.Cells.Item(5, 8).Value = "MSD Central America & The Caribbean"

$ws.Cells.Item(6, 1).Value = "1327240"
$ws.Cells.Item(6, 2).Value = "https://aiesec.org/opportunity/global-talent/1327240"
$ws.Cells.Item(6, 3).Value = "Business Development"
$ws.Cells.Item(6, 4).Value = "Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt"
$ws.Cells.Item(6, 5).Value = "No"
$ws.Cells.Item(6, 6).Value = "5 applicants"
$ws.Cells.Item(6, 7).Value = "9 - 12 Weeks"
$ws.Cells.Item(6, 8).Value = "Vista"

$ws.Cells.Item(7, 1).Value = "1326944"
$ws.Cells.Item(7, 2).Value = "https://aiesec.org/opportunity/global-talent/1326944"
$ws.Cells.Item(7, 3).Value = "Digital Marketing Executive"
$ws.Cells.Item(7, 4).Value = "Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt"
$ws.Cells.Item(7, 5).Value = "No"
$ws.Cells.Item(7, 6).Value = "13 applicants"
$ws.Cells.Item(7, 7).Value = "9 - 12 Weeks"
$ws.Cells.Item(7, 8).Value = "TIE innovated Solution"

$ws.Cells.Item(8, 1).Value = "1326481"
$ws.Cells.Item(8, 2).Value = "https://aiesec.org/opportunity/global-talent/1326481"
$ws.Cells.Item(8, 3).Value = "Global Duty Billing Data Analytics Expert"
$ws.Cells.Item(8, 4).Value = "Maastricht, Netherlands"
$ws.Cells.Item(8, 5).Value = "Yes"
$ws.Cells.Item(8, 6).Value = "127 applicants"
$ws.Cells.Item(8, 7).Value = "6 - 18 Months"
$ws.Cells.Item(8, 8).Value = "DHL Group"
$ws.Cells.Item(8, 5).Interior.Color = 65535

$ws.Cells.Item(9, 1).Value = "1326043"
$ws.Cells.Item(9, 2).Value = "https://aiesec.org/opportunity/global-talent/1326043"
$ws.Cells.Item(9, 3).Value = "ACE Program | Portuguese Financial Analyst"
$ws.Cells.Item(9, 4).Value = "Kolkata, West Bengal, India"
$ws.Cells.Item(9, 5).Value = "Yes"
$ws.Cells.Item(9, 6).Value = "24 applicants"
$ws.Cells.Item(9, 7).Value = "6 - 18 Months"
$ws.Cells.Item(9, 8).Value = "Tata Consultancy Services Ltd."
$ws.Cells.Item(9, 5).Interior.Color = 65535

$ws.Cells.Item(10, 1).Value = "1326042"
$ws.Cells.Item(10, 2).Value = "https://aiesec.org/opportunity/global-talent/1326042"
$ws.Cells.Item(10, 3).Value = "ACE Program | Portuguese Financial Analyst"
$ws.Cells.Item(10, 4).Value = "Thane, Maharashtra, India"
$ws.Cells.Item(10, 5).Value = "Yes"
$ws.Cells.Item(10, 6).Value = "22 applicants"
$ws.Cells.Item(10, 7).Value = "6 - 18 Months"
$ws.Cells.Item(10, 8).Value = "Tata Consultancy Services Ltd."
$ws.Cells.Item(10, 5).Interior.Color = 65535

$ws.Cells.Item(11, 1).Value = "1326041"
$ws.Cells.Item(11, 2).Value = "https://aiesec.org/opportunity/global-talent/1326041"
$ws.Cells.Item(11, 3).Value = "ACE Program | Spanish Financial Analyst"
$ws.Cells.Item(11, 4).Value = "Thane, Maharashtra, India"
$ws.Cells.Item(11, 5).Value = "Yes"
$ws.Cells.Item(11, 6).Value = "18 applicants"
$ws.Cells.Item(11, 7).Value = "6 - 18 Months"
$ws.Cells.Item(11, 8).Value = "Tata Consultancy Services Ltd."
$ws.Cells.Item(11, 5).Interior.Color = 65535

$ws.Cells.Item(12, 1).Value = "1312732"
$ws.Cells.Item(12, 2).Value = "https://aiesec.org/opportunity/global-talent/1312732"
$ws.Cells.Item(12, 3).Value = "Marketing Executive"
$ws.Cells.Item(12, 4).Value = "Hyderabad, Telangana, India"
$ws.Cells.Item(12, 5).Value = "No"
$ws.Cells.Item(12, 6).Value = "5 applicants"
$ws.Cells.Item(12, 7).Value = "9 - 12 Weeks"
$ws.Cells.Item(12, 8).Value = "TERICSOFT TECHNOLOGY SOLUTIONS PVT. LTD."

$ws.Cells.Item(13, 1).Value = "1309418"
$ws.Cells.Item(13, 2).Value = "https://aiesec.org/opportunity/global-talent/1309418"
$ws.Cells.Item(13, 3).Value = "International Marketing Intern"
$ws.Cells.Item(13, 4).Value = "Luzhu District, Kaohsiung City, Taiwan 821"
$ws.Cells.Item(13, 5).Value = "No"
$ws.Cells.Item(13, 6).Value = "95 applicants"
$ws.Cells.Item(13, 7).Value = "6 - 18 Months"
$ws.Cells.Item(13, 8).Value = "Green Power Engineering Corporation"
$ws.Columns.Item(3).ColumnWidth = 44.166666666666664
$ws.Columns.Item(4).ColumnWidth = 56.166666666666664
